$d = $word.ActiveDocument

# 1. "...and social and environment impact are all topics..."
#    -> "...and social and environmental impacts are all topics..."
$d.Content.Find.Execute(
    " and social and environment impact are all topics that are important to consider when contemplating how we manage the internet",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and social and environmental impacts are all topics that are important to consider when contemplating how we manage the internet",
    2) | Out-Null

# 2. "...in contact with a cloud..." -> "...in contact with the Cloud..."
$d.Content.Find.Execute(
    " In fact, unless specifically designed not to, almost all current data is in contact with a cloud " + [char]0x2013 + " in some form or another.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " In fact, unless specifically designed not to, almost all current data is in contact with the Cloud " + [char]0x2013 + " in some form or another.",
    2) | Out-Null

# 3. "...offering consumers cloud-based services" -> "...offering consumers high quality cloud-based services"
$d.Content.Find.Execute(
    "Companies are now putting a greater focus on offering consumers cloud-based services",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Companies are now putting a greater focus on offering consumers high quality cloud-based services",
    2) | Out-Null

# 4. "...redundancy options and unparalleled..." -> "...redundancy and unparalleled..."
$d.Content.Find.Execute(
    "thorough backing up and redundancy options and unparalleled",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "thorough backing up and redundancy and unparalleled",
    2) | Out-Null

# 5. "...over long distances." -> "...over long distances using the Internet Protocol."
$d.Content.Find.Execute(
    [char]0x2013 + " a system that allows users to share data between devices over long distances.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x2013 + " a system that allows users to share data between devices over long distances using the Internet Protocol.",
    2) | Out-Null

# 6. "...workplace or school setting." -> "...vocational or educational setting."
$d.Content.Find.Execute(
    "groups of people in a professional capacity. Specifically focuses on productivity in a workplace or school setting.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "groups of people in a professional capacity. Specifically focuses on productivity in a vocational or educational setting.",
    2) | Out-Null

# 7. "Examples include Amazon Web Services..." -> "Examples include Microsoft Azure, Amazon Web Services..."
$d.Content.Find.Execute(
    "PaaS or suites that allow development, management and deployment of web and mobile applications. Examples include Amazon Web Services and Google App Engine.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PaaS or suites that allow development, management and deployment of web and mobile applications. Examples include Microsoft Azure, Amazon Web Services and Google App Engine.",
    2) | Out-Null

# 8. "...customer data on pre-structured data pools..." -> "...customer data from pre-structured data pools..."
$d.Content.Find.Execute(
    "market analytics on massive amounts of customer data on pre-structured data pools like Google",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "market analytics on massive amounts of customer data from pre-structured data pools like Google",
    2) | Out-Null

# 9. Big rewrite: "...ever before but is there a cost to this much benefit? ... Well, in ways there is and"
#    -> "...ever before. But is there a cost to such a boon? ... Well, certainly in some ways there is an"
$d.Content.Find.Execute(
    "and there are more software development options than ever before but is there a cost to this much benefit? And where are the nay-sayers, if there are any? Well, in ways there is and",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and there are more software development options than ever before. But is there a cost to such a boon? And where are the nay-sayers, if there are any? Well, certainly in some ways there is an",
    2) | Out-Null

# 10. "...super-famous celebrities nude photos..." -> "...super-famous celebrities' nude photos..."
$d.Content.Find.Execute(
    "thousands of super-famous celebrities nude photos were leaked to the public",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "thousands of super-famous celebrities" + [char]0x2019 + " nude photos were leaked to the public",
    2) | Out-Null

Write-Host "All replacements applied"
